$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 contains a single user record; update the randomized test data values
$ws.Range("A2").Value = "QzBON561"
$ws.Range("B2").Value = 23072734
$ws.Range("C2").Value = "gtbcicg49"
$ws.Range("D2").Value = "Rf5&4Fp%"
$ws.Range("F2").Value = "EdwytpvL"
$ws.Range("G2").Value = "EIzt"
